{"js": "// The scraped agenda page was a single run of run-on text. Rebuild it as\n// a structured paragraph: one run holding many text segments separated by\n// manual line breaks (\"\\v\" inserts a <w:br/> between <w:t> runs, the same\n// way Word represents Shift+Enter), followed by a second paragraph holding\n// the \"End of Page 1\" marker (also wrapped in line breaks).\nconst body = context.document.body;\nbody.clear();\n\nconst lines = [\n  \"cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4445\",\n  \"6/23/23, 5:41 PM\",\n  \"ORLEANS\",\n  \"NEW O\",\n  \"CITY O\",\n  \"LOUISIANA\",\n  \"CRIMINAL JUSTICE COMMITTEE\",\n  \"MEETING AGENDA\",\n  \"TUESDAY, MARCH 7, 2023\",\n  \"10:00 A.M.\",\n  \"1. Roll Call\",\n  \"Cal No .: 34,062 - By Councilmember Thomas (By Request) - To Present: A\",\n  \"4.\",\n  \"Representative of the New Orleans Health Department- For Discussion and Vote: An\",\n  \"Ordinance to amend and reordain Section 54-28 and Section 54-338 of the Code of the\",\n  \"City of New Orleans, to allow for custodial arrest for adults for the careless storage of a\",\n  \"firearm that results in a minor either causing, or attempting to cause, the injury or\",\n  \"death of themselves or any other person, or causes the firearm to fire; and provides for\",\n  \"a penalty for any person who fails to secure a firearm that a minor obtains, which\",\n  \"results in a minor either causing, or attempting to cause, the injury or death of\",\n  \"themselves or any other person or causes the firearm to fire.\",\n  \"T Ord. Cal. No. 34,062\",\n  \"2. Approval of the minutes of the February 7, 2023 meeting\",\n  \"Motion M-23-I - By Councilmember Thomas (By Request) - To Present: Courtney\",\n  \"3.\",\n  \"Story, CAO's office - For Discussion and Vote: Motion to approve a proposed\",\n  \"Amendment No. 2 of a professional services agreement between the City of New\",\n  \"Orleans and American Traffic Solutions, Inc. to assist the New Orleans Police\",\n  \"Department by providing traffic safety program management services to the City, for a\",\n  \"nineteen (19) months period, with a total compensation amount not to exceed\",\n  \"$11,300,000.00.\",\n  \"TM-23-1\",\n  \"5. Adjournment\",\n  \"Public Comment\",\n  \"1/1\",\n  \"https://cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4445\"\n];\n\nconst firstParagraph = body.paragraphs.getFirst();\nfirstParagraph.insertText(lines.join(\"\\v\"), \"Start\");\nawait context.sync();\n\nconst secondParagraph = firstParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\nsecondParagraph.insertText(\"\\v\" + \"---------- End of Page 1 ----------\" + \"\\v\", \"Start\");\nawait context.sync();\n", "ps1": "# Rebuild the scraped agenda text as a structured, line-broken run\n# (one paragraph holding many <w:t> segments joined by line breaks),\n# followed by a second paragraph containing the \"End of Page 1\" marker.\n$d = $word.ActiveDocument\n\n$lineBreak = [char]11\n\n$lines = @(\n    'cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4445',\n    '6/23/23, 5:41 PM',\n    'ORLEANS',\n    'NEW O',\n    'CITY O',\n    'LOUISIANA',\n    'CRIMINAL JUSTICE COMMITTEE',\n    'MEETING AGENDA',\n    'TUESDAY, MARCH 7, 2023',\n    '10:00 A.M.',\n    '1. Roll Call',\n    'Cal No .: 34,062 - By Councilmember Thomas (By Request) - To Present: A',\n    '4.',\n    'Representative of the New Orleans Health Department- For Discussion and Vote: An',\n    'Ordinance to amend and reordain Section 54-28 and Section 54-338 of the Code of the',\n    'City of New Orleans, to allow for custodial arrest for adults for the careless storage of a',\n    'firearm that results in a minor either causing, or attempting to cause, the injury or',\n    'death of themselves or any other person, or causes the firearm to fire; and provides for',\n    'a penalty for any person who fails to secure a firearm that a minor obtains, which',\n    'results in a minor either causing, or attempting to cause, the injury or death of',\n    'themselves or any other person or causes the firearm to fire.',\n    'T Ord. Cal. No. 34,062',\n    '2. Approval of the minutes of the February 7, 2023 meeting',\n    'Motion M-23-I - By Councilmember Thomas (By Request) - To Present: Courtney',\n    '3.',\n    'Story, CAO''s office - For Discussion and Vote: Motion to approve a proposed',\n    'Amendment No. 2 of a professional services agreement between the City of New',\n    'Orleans and American Traffic Solutions, Inc. to assist the New Orleans Police',\n    'Department by providing traffic safety program management services to the City, for a',\n    'nineteen (19) months period, with a total compensation amount not to exceed',\n    '$11,300,000.00.',\n    'TM-23-1',\n    '5. Adjournment',\n    'Public Comment',\n    '1/1',\n    'https://cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4445'\n)\n\n$d.Content.Text = ($lines -join $lineBreak)\n\n$secondParagraph = $d.Paragraphs.Add()\n$secondParagraph.Range.Text = $lineBreak + '---------- End of Page 1 ----------' + $lineBreak\n\n"}
